$wb = $excel.ActiveWorkbook

# 1. Rename "safe_column_name_test" -> "safe_header_name_test"
$ws2 = $wb.Worksheets.Item("safe_column_name_test")
$ws2.Name = "safe_header_name_test"

# 2. Add a new sheet "date_test" after the last sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "date_test"

# 3. Populate the "date_test" sheet
$ws3.Range("A1").Value = "date"
$ws3.Range("B1").Value = "plaincol"

$ws3.Range("A2").Value = 29580
$ws3.Range("A2").NumberFormat = "m/d/yy"
$ws3.Range("B2").Value = "it will still parse the dates below as date even if plaincol is not in the default --dates-whitelist because the cell format was set to date"
$ws3.Range("B2").WrapText = $true

$ws3.Range("A3").Value = 37145.354166666664
$ws3.Range("A3").NumberFormat = '[$-409]m/d/yy\ h:mm\ AM/PM;@'
$ws3.Range("B3").Value = 37145
$ws3.Range("B3").NumberFormat = "m/d/yy"

$ws3.Range("A4").Value = "not a date"
$ws3.Range("B4").Value = 37145.354166666664
$ws3.Range("B4").NumberFormat = "m/d/yy h:mm"

$ws3.Range("A5").Value = "Wednesday, Mar-14-2012"
$ws3.Range("A5").NumberFormat = '[$-F800]dddd\,\ mmmm\ dd\,\ yyyy'
$ws3.Range("B5").Value = "the date below is not parsed as a date coz we didn't explicitly set the cell format to a date format and ""plaincol"" is not in the --dates-whitelist"
$ws3.Range("B5").WrapText = $true

$ws3.Range("A6").Value = 37145
$ws3.Range("A6").NumberFormat = '[$-F800]dddd\,\ mmmm\ dd\,\ yyyy'
$ws3.Range("B6").Value = "9/11/01 8:30 am"
$ws3.Range("B6").NumberFormat = "@"

# column widths / row heights
$ws3.Columns.Item(1).ColumnWidth = 26.83
$ws3.Columns.Item(2).ColumnWidth = 41.5
$ws3.Rows.Item(2).RowHeight = 58
$ws3.Rows.Item(5).RowHeight = 48

Write-Host "Sheets:"
foreach ($s in $wb.Worksheets) {
    Write-Host $s.Name
}
